# chemcurier v.0.9.4 v.1.6 with download at 04/11/2023
# Adds Tyre Size / Model / Param / Sales value / Date_of_sales / Contragent
# columns (E:J) to the tyres sheet, and drops the stray empty C1 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The lone placeholder cell in C1 is removed entirely (row 1 keeps only A1/B1,
# plus the new header cells below).
$ws.Range("C1").ClearContents()

# New header row (E1:J1)
$ws.Range("E1").Value = "Tyre Size"
$ws.Range("F1").Value = "Model"
$ws.Range("G1").Value = "Param"
$ws.Range("H1").Value = "Sales value"
$ws.Range("I1").Value = "Date_of_sales"
$ws.Range("J1").Value = "Contragent"

# Data rows 2-21: E = tyre size, F = model, G = comma-joined param list,
# H = sales value, I = date of sales, J = contragent.
$data = @(
    @(2, "35/65-33", "ФБел-283", "42, 30, груз, сер"),
    @(3, "205/55R16", "BEL-262", "сер, легк, б/к"),
    @(4, "205/55R16", "BEL-317", "сер, легк, б/к"),
    @(5, "205/55R16", "BEL-317S", "ошип, сер"),
    @(6, "235/75R15", "BEL-1001", "сер, легк"),
    @(7, "155/65R13", "BEL-1002", "сер, легк"),
    @(8, "205/55R16", "BEL-1004", "сер, легк"),
    @(9, "225/50R17", "BEL-1005", "сер, легк"),
    @(10, "24.00R35", "Бел-202", "210B, H, Type, C, сер"),
    @(11, "24.00R35", "Бел-212", "груз, Type, LS-2, сер"),
    @(12, "21.00R35", "Бел-200", "202B, Type, C, сер"),
    @(13, "21.00R35", "Бел-210", "202B, H, Type, LS-2, C, сер"),
    @(14, "14.00R20", "BEL-248", "груз, сер, б/к"),
    @(15, "175/70R13", "Бел-103", "сер, легк, б/к"),
    @(16, "175/70R13", "Бел-100", "сер, легк, б/к"),
    @(17, "195/65R15", "Бел-119", "сер, легк"),
    @(18, "210/80R16", "Бел-777", "сер, легк"),
    @(19, "215/65R16C", "Бел-1000", "сер, легк"),
    @(20, "205/55R16", "Бел-1001", "сер, легк"),
    @(21, "225/50R17", "Бел-1005", "сер, легк")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 5).Value = $item[1]
    $ws.Cells.Item($r, 6).Value = $item[2]
    $ws.Cells.Item($r, 7).Value = $item[3]
    $ws.Cells.Item($r, 8).Value = 2
    $dateCell = $ws.Cells.Item($r, 9)
    $dateCell.Value = 45138
    $dateCell.NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r, 10).Value = "нет данных"
}
